$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20 continues the existing table; copy the formatting of the
# date cell above it (A19) onto A20 before writing the new values.
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -2.06674933094535
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -0.3099928749133896
